# Update stats for 2026-01 (row 26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6478
$ws.Range("D26").Value = 6041261
$ws.Range("E26").Value = 932.5811979005866
$ws.Range("F26").Value = 9.481155991211754
$ws.Range("H26").Value = 25.80946385326015
